{"js": "// \"vide\u00f3 sorrend a scriptbe\"\n// Append the topology walk-through paragraph, a \"Szerverek sz\u00f6veg:\" heading\n// paragraph, and a trailing blank paragraph at the end of the document body\n// (right before the section properties), mirroring the diff that adds three\n// new <w:p> blocks after the final existing paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst topologyText =\n  \"Topologia bemutat\u00e1sa: laguna -> k\u00e1v\u00e9h\u00e1z -> new york (vide\u00f3k r\u00f3la) -> home office -> routerek (\u00e9s isp) -> \u00fczemeltet\u0151k -> szerverek \u00e9s szolg\u00e1ltat\u00e1sok\";\n\nconst topologyParagraph = lastParagraph.insertParagraph(topologyText, \"After\");\nconst serverLabelParagraph = topologyParagraph.insertParagraph(\"Szerverek sz\u00f6veg:\", \"After\");\n// Trailing empty paragraph, matching the diff's final <w:p/>.\nserverLabelParagraph.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# \"vide\u00f3 sorrend a scriptbe\"\n# Append the topology walk-through paragraph, a \"Szerverek sz\u00f6veg:\" heading\n# paragraph, and a trailing blank paragraph at the end of the document body\n# (right before the section properties), mirroring the diff that adds three\n# new <w:p> blocks after the final existing paragraph.\n\n$d = $word.ActiveDocument\n\n$topologyText = \"Topologia bemutat\u00e1sa: laguna -> k\u00e1v\u00e9h\u00e1z -> new york (vide\u00f3k r\u00f3la) -> home office -> routerek (\u00e9s isp) -> \u00fczemeltet\u0151k -> szerverek \u00e9s szolg\u00e1ltat\u00e1sok\"\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n$topologyParagraph = $d.Paragraphs.Last\n$topologyParagraph.Range.Text = $topologyText\n\n$topologyParagraph.Range.InsertParagraphAfter()\n$serverLabelParagraph = $d.Paragraphs.Last\n$serverLabelParagraph.Range.Text = \"Szerverek sz\u00f6veg:\"\n\n# Trailing empty paragraph, matching the diff's final <w:p/>.\n$serverLabelParagraph.Range.InsertParagraphAfter()\n"}
